# Insert a new row above current row 3, shifting existing rows 3-13 down to 4-14.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Insert()

# Fill the new row 3 with its data.
# Column A keeps the running numeric sequence (row3 -> 2).
$ws.Range("A3").Value = 2

# Add the new shared strings in the same order they appear in the target
# workbook (조감도 = index 18, 지도 = index 19): set D before C.
$ws.Range("D3").Value = "조감도"
$ws.Range("C3").Value = "지도"

# Renumber column A for the rows that were pushed down (now rows 4-14),
# so the sequence stays 1..13.
for ($i = 4; $i -le 14; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

# Update the active selection to match the target workbook.
$ws.Range("F7").Select()
